$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.269.24"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.907.91"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'307.78"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5297"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").Value = "'0.3820"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "'0.07293"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'22.08"
$ws.Range("E10").Value = "  +4.40%  "
$ws.Range("D11").Value = "'0.9018"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'0.08193"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "'95.73"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "'5.345"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "'1.003"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'0.000008653"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'14.79"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "27.308.36"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "1.195.82"
$ws.Range("E20").Value = "  -37.05%  "
$ws.Range("D21").Value = "'5.056"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").Value = "'6.517"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "'149.99"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").Value = "'18.24"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'1.745"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'116.97"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "'4.819"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "'4.810"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").Value = "'0.09288"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'0.8362"
$ws.Range("E32").Value = "  +4.88%  "
$ws.Range("D33").Value = "'0.05067"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'1.224"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("D35").Value = "'3.008"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("D36").Value = "'3.359"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("D38").Value = "'0.5741"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'0.02006"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "'1.077"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "'9.307"
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("D42").Value = "'6.560"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "'117.25"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "'0.1524"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "'0.4924"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "'10.14"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").Value = "'1.636"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "'38.71"
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("D50").Value = "'0.06159"
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("D51").Value = "'63.61"
$ws.Range("E51").Value = "  -0.29%  "
